$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quantities in the I/O list to reflect the reduced scope of the
# eplan schematic start (16 instead of 32 mass valves, 23 outputs instead of
# 39, and 3 instead of 4 pressure transducers).
$ws.Range("C5").Value = "16x válvulas das massas"
$ws.Range("C1").Value = "Saídas - 23"
$ws.Range("A13").Value = "3x transdutor de pressão"

# Leave the selection where the last edit was made.
$ws.Range("A14").Select()
